$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '76.527.53'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.106.51'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +5.06%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '199.43'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '623.50'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.62%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +6.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.553'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.466'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.26'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.02%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.651.62'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '29.60'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.68%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000201'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +5.50%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '76.389.05'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.072.42'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +4.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.64'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.18'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.87%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.75'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +20.30%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '387.32'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.55'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.79%  '
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.61'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +6.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.246.02'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.94%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '72.80'
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.39'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +6.48%  '
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('E29').Value = '  +3.42%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.44'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.46'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.63%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '512.49'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.95'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +6.72%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.133'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +19.29%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '21.02'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '163.74'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '196.54'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +8.64%  '
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.380'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.103'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.86%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.32'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +7.53%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.803'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +20.30%  '
$ws.Range('E46').Value = '  +7.44%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.70'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.51'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +7.66%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '41.16'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.84%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.607'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.45%  '
$ws.Range('E51').Value = '  +1.58%  '
